$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed the new lookup/category values first (matches the order they were
# first introduced into the workbook's shared-string table).
$ws.Range("A2").Value = "REFRIGERATOR"
$ws.Range("A4").Value = "FREEZER"
$ws.Range("B7").Value = "ELECTROLUX"
$ws.Range("B2").Value = "SANDEN INTERCOOL"
$ws.Range("B5").Value = "KOKUSAN"
$ws.Range("C5").Value = "H-19α"

# Row 2: REFRIGERATOR / SANDEN INTERCOOL / YPM-165P / YPM165PS-160300320 / LAB-001 / 2025-06-11 / CH25062842
$ws.Range("A2").Value = "REFRIGERATOR"
$ws.Range("B2").Value = "SANDEN INTERCOOL"
$ws.Range("C2").Value = "YPM-165P"
$ws.Range("D2").Value = "YPM165PS-160300320"
$ws.Range("E2").Value = "LAB-001"
$ws.Range("G2").Value = "CH25062842"

# Row 3: REFRIGERATOR / SANDEN INTERCOOL / OEM-1205I / OEM1205I-221001861 / LAB-002 / 2025-06-11 / CH25062843
$ws.Range("A3").Value = "REFRIGERATOR"
$ws.Range("B3").Value = "SANDEN INTERCOOL"
$ws.Range("C3").Value = "OEM-1205I"
$ws.Range("D3").Value = "OEM1205I-221001861"
$ws.Range("E3").Value = "LAB-002"
$ws.Range("G3").Value = "CH25062843"

# Row 4: FREEZER / FRESHER / FF-182WD / FF182WD-210900016 / LAB-003 / 2025-06-11 / CH25062844
$ws.Range("A4").Value = "FREEZER"
$ws.Range("B4").Value = "FRESHER"
$ws.Range("C4").Value = "FF-182WD"
$ws.Range("D4").Value = "FF182WD-210900016"
$ws.Range("E4").Value = "LAB-003"
$ws.Range("G4").Value = "CH25062844"

# Row 5: CENTRIFUGE / KOKUSAN / H-19α / 147925 / LAB-004 / 2025-06-11 / CF25062845
$ws.Range("A5").Value = "CENTRIFUGE"
$ws.Range("B5").Value = "KOKUSAN"
$ws.Range("C5").Value = "H-19α"
$ws.Range("D5").Value = 147925
$ws.Range("E5").Value = "LAB-004"
$ws.Range("G5").Value = "CF25062845"

# Row 6: CENTRIFUGE / KOKUSAN / H-19α / 150753 / LAB-005 / 2025-06-11 / CF25062846
$ws.Range("A6").Value = "CENTRIFUGE"
$ws.Range("B6").Value = "KOKUSAN"
$ws.Range("C6").Value = "H-19α"
$ws.Range("D6").Value = 150753
$ws.Range("E6").Value = "LAB-005"
$ws.Range("G6").Value = "CF25062846"

# Row 7: REFRIGERATOR / ELECTROLUX / EUM0930AD-TH / 41361175 / LAB-026 / 2025-06-11 / CH25062847
$ws.Range("A7").Value = "REFRIGERATOR"
$ws.Range("B7").Value = "ELECTROLUX"
$ws.Range("C7").Value = "EUM0930AD-TH"
$ws.Range("D7").Value = 41361175
$ws.Range("E7").Value = "LAB-026"
$ws.Range("G7").Value = "CH25062847"

# Rows 8-11: clear out the previously populated data (now blank rows, same as row 12+)
$ws.Range("A8:G11").ClearContents()

# Calibration Date column (F) for rows 2-7 is 2025-06-11 (Excel serial 45819)
$ws.Range("F2:F7").Value = 45819

# Restore the saved selection to F12
$ws.Range("F12").Select()
